$wb = $excel.ActiveWorkbook

# --- Sheet: "Transfer Time (s)" ---
$ws1 = $wb.Worksheets.Item("Transfer Time (s)")
$ws1.Range("D4").Value = 0.07475761151313781
$ws1.Range("E4").Value = 0.02206291961922436
$ws1.Range("D5").Value = 0.08424967050552368
$ws1.Range("E5").Value = 0.01959415137762406
$ws1.Range("D6").Value = 0.09851876497268677
$ws1.Range("E6").Value = 0.03453666244963857
$ws1.Range("D7").Value = 0.3237078189849854

# --- Sheet: "Throughput (bps)" ---
$ws2 = $wb.Worksheets.Item("Throughput (bps)")
$ws2.Range("D4").Value = 143886.396252974
$ws2.Range("E4").Value = 25727.35812289039
$ws2.Range("D5").Value = 1288859.973798563
$ws2.Range("E5").Value = 250854.3482564982
$ws2.Range("D6").Value = 11844353.21478873
$ws2.Range("E6").Value = 2833977.45282994
$ws2.Range("D7").Value = 36648994.14349551
